$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text could be misread as numbers/dates by Excel's auto-detection
# are given an explicit Text number format (quote-prefix behavior) before the
# value is assigned, so they are stored as text -- matching the source data which
# is all plain text (inlineStr) in the workbook.

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

$ws.Range("D2").Value = "29.876.14"
$ws.Range("E2").Value = "  +0.16%  "
$ws.Range("D3").Value = "1.887.70"
$ws.Range("E3").Value = "  -0.22%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "0.7683"
$ws.Range("E5").Value = "  -0.79%  "
$ws.Range("D6").Value = "242.72"
$ws.Range("E6").Value = "  -0.78%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").Value = "0.3126"
$ws.Range("E8").Value = "  -0.61%  "
$ws.Range("D9").Value = "25.63"
$ws.Range("E9").Value = "  +1.05%  "
$ws.Range("D10").Value = "0.07171"
$ws.Range("E10").Value = "  -3.88%  "
$ws.Range("D11").Value = "0.08579"
$ws.Range("E11").Value = "  +5.69%  "
$ws.Range("E12").Value = "  -0.48%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.905.23"
$ws.Range("E13").Value = "  +0.35%  "
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value = "5.364"
$ws.Range("E14").Value = "  -1.64%  "
$ws.Range("D15").Value = "93.62"
$ws.Range("E15").Value = "  +1.54%  "
$ws.Range("D16").Value = "6.146"
$ws.Range("E16").Value = "  -0.96%  "
$ws.Range("D17").Value = "29.841.90"
$ws.Range("E17").Value = "  +0.06%  "
$ws.Range("E18").Value = "  -1.47%  "
$ws.Range("D19").Value = "244.35"
$ws.Range("E19").Value = "  +0.12%  "
$ws.Range("D20").Value = "0.000007805"
$ws.Range("E20").Value = "  -0.94%  "
$ws.Range("D21").Value = "2.142.02"
$ws.Range("E21").Value = "  +0.82%  "
$ws.Range("E22").Value = "  -0.11%  "
$ws.Range("D23").Value = "8.007"
$ws.Range("E23").Value = "  -1.03%  "
$ws.Range("E24").Value = "  -0.03%  "
$ws.Range("D25").Value = "0.1636"
$ws.Range("E25").Value = "  +3.39%  "
$ws.Range("E26").Value = "  -0.50%  "
$ws.Range("D27").Value = "162.96"
$ws.Range("E27").Value = "  +0.14%  "
$ws.Range("D28").Value = "18.72"
$ws.Range("E28").Value = "  -0.45%  "
$ws.Range("D29").Value = "2.032"
$ws.Range("E29").Value = "  -0.58%  "
$ws.Range("D30").Value = "1.469"
$ws.Range("E30").Value = "  +2.63%  "
$ws.Range("D31").Value = "1.534"
$ws.Range("E31").Value = "  -1.07%  "
$ws.Range("E32").Value = "  +0.54%  "
$ws.Range("D33").Value = "4.093"
$ws.Range("E33").Value = "  -0.11%  "
$ws.Range("D34").Value = "0.05450"
$ws.Range("E34").Value = "  -1.08%  "
$ws.Range("D35").Value = "1.240"
$ws.Range("E35").Value = "  -1.11%  "
$ws.Range("D36").Value = "0.7424"
$ws.Range("E36").Value = "  -1.79%  "
$ws.Range("D37").Value = "1.000"
$ws.Range("E37").Value = "  +0.09%  "
$ws.Range("E38").Value = "  +2.12%  "
$ws.Range("D39").Value = "0.01953"
$ws.Range("E39").Value = "  +1.45%  "
$ws.Range("D40").Value = "2.785"
$ws.Range("E40").Value = "  -0.14%  "
$ws.Range("D41").Value = "0.4469"
$ws.Range("E41").Value = "  +0.45%  "
$ws.Range("D42").Value = "1.106.31"
$ws.Range("E42").Value = "  -4.87%  "
$ws.Range("D43").Value = "73.12"
$ws.Range("E43").Value = "  -1.03%  "
$ws.Range("D44").Value = "6.068"
$ws.Range("E44").Value = "  +1.89%  "
$ws.Range("D45").Value = "0.8521"
$ws.Range("E45").Value = "  +0.44%  "
$ws.Range("D46").Value = "1.001"
$ws.Range("E46").Value = "  -0.01%  "
$ws.Range("D47").Value = "102.57"
$ws.Range("E47").Value = "  +0.47%  "
$ws.Range("D48").Value = "7.655"
$ws.Range("E48").Value = "  +1.50%  "
$ws.Range("D49").Value = "1.862"
$ws.Range("E49").Value = "  -2.06%  "
$ws.Range("D50").Value = "3.011"
$ws.Range("E50").Value = "  -2.32%  "
$ws.Range("D51").Value = "2.039.79"
$ws.Range("E51").Value = "  +0.49%  "
